# Weekly work breakdown - fill in Week 1 (row 4) tasks for each team member column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Lit Review"
$ws.Range("D4").Value = "Aims and objectives + Class and sequance diagrams"
$ws.Range("E4").Value = "Sidebar GUI"
$ws.Range("F4").Value = "Filtering"
$ws.Range("G4").Value = "Sorting"
$ws.Range("H4").Value = "Machine learning"
$ws.Range("I4").Value = "Implementation report"
$ws.Range("J4").Value = "Implementation report"

# Widen the columns that now hold longer text (Description and Filtering columns)
$ws.Columns.Item(4).ColumnWidth = 52.0
$ws.Columns.Item(6).ColumnWidth = 23.5

# Update the current cell selection to where editing finished
$ws.Range("I6").Select() | Out-Null

# Set print page setup (paper size / orientation) as configured by the author
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
